$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (cycle determining function
# now accounts for trailing zero-current bpoints more than an hour apart).
$ws.Range("B2").Value = 0.60409275
$ws.Range("C2").Value = 0.20266954
$ws.Range("D2").Value = 0.12120013
$ws.Range("E2").Value = 0.07203758

$ws.Range("B3").Value = 0.23597553
$ws.Range("C3").Value = 0.40452926
$ws.Range("D3").Value = 0.30179174
$ws.Range("E3").Value = 0.05770347

$ws.Range("B4").Value = 0.01916627
$ws.Range("C4").Value = 0.0871682
$ws.Range("D4").Value = 0.17918482
$ws.Range("E4").Value = 0.71448071

$ws.Range("B5").Value = 0.14076545
$ws.Range("C5").Value = 0.305633
$ws.Range("D5").Value = 0.39782331
$ws.Range("E5").Value = 0.15577824
